$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D keeps its text (inline-string) representation instead of
# being auto-converted to a number by Excel when we assign numeric-looking text.
$ws.Range("D2:D11").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Move Robot15 to location (8, 6) and remove the screws."
$ws.Range("B2").Value = 22.995794
$ws.Range("C2").Value = 3855
$ws.Range("D2").Value = "0.0081"
$ws.Range("E2").Value = "fab7db7a-f521-4140-8eb4-0f0c395ac018"

# Row 3
$ws.Range("A3").Value = "Move to location (10, 2) and remove the toolkit."
$ws.Range("B3").Value = 25.411137
$ws.Range("C3").Value = 3953
$ws.Range("D3").Value = "0.00891"
$ws.Range("E3").Value = "288e9003-0b70-446c-9a60-1a6c03025922"

# Row 4
$ws.Range("A4").Value = "Move to location (6, 6) and remove the liquid spill."
$ws.Range("B4").Value = 24.897036
$ws.Range("C4").Value = 3798
$ws.Range("D4").Value = "0.00753"
$ws.Range("E4").Value = "766135a4-8835-4689-91e1-b9a16aae055e"

# Row 5
$ws.Range("A5").Value = "Move to location (5, 3) and remove the large debris."
$ws.Range("B5").Value = 23.161379
$ws.Range("C5").Value = 3919
$ws.Range("D5").Value = "0.00804"
$ws.Range("E5").Value = "94a75475-33cb-4a6b-9b7b-466174d99cb5"

# Row 6
$ws.Range("A6").Value = "Move to location (6, 10) and remove the dust."
$ws.Range("B6").Value = 31.108275
$ws.Range("C6").Value = 3875
$ws.Range("D6").Value = "0.00846"
$ws.Range("E6").Value = "8574efd5-4df5-44a2-b70b-34572fb96b50"

# Row 7
$ws.Range("A7").Value = "Move to location (4, 8) and remove the grass."
$ws.Range("B7").Value = 30.749819
$ws.Range("C7").Value = 3888
$ws.Range("D7").Value = "0.00822"
$ws.Range("E7").Value = "29a8b512-8be1-4602-8b40-2e0a8e63f366"

# Row 8
$ws.Range("A8").Value = "Move to location (8, 7) and remove the small debris."
$ws.Range("B8").Value = 31.511673
$ws.Range("C8").Value = 3919
$ws.Range("D8").Value = "0.00828"
$ws.Range("E8").Value = "05435f9d-2dd5-4042-8c1d-5dd0c1f33e59"

# Row 9
$ws.Range("A9").Value = "Move to location (1, 10) and remove the vehicle."
$ws.Range("B9").Value = 31.595244
$ws.Range("C9").Value = 4512
$ws.Range("D9").Value = "0.00915"
$ws.Range("E9").Value = "09a1c6e2-cd69-4b43-9641-54831cd4bb9c"

# Row 10
$ws.Range("A10").Value = "Move to location (2, 12) and remove the construction materials."
$ws.Range("B10").Value = 30.238777
$ws.Range("C10").Value = 3890
$ws.Range("D10").Value = "0.00795"
$ws.Range("E10").Value = "2b52dc19-88f3-4e0c-84b1-c2550a97f1be"

# Row 11
$ws.Range("A11").Value = "Move to location (8, 9) and remove the tree branches."
$ws.Range("B11").Value = 30.231626
$ws.Range("C11").Value = 3879
$ws.Range("D11").Value = "0.00762"
$ws.Range("E11").Value = "d2772432-8d5a-4185-aa51-12938fbda2f3"
